$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - strict mode violation error (3 elements), new time
$f2 = @'
Locator.click: Error: strict mode violation: locator("button:has-text('Next')") resolved to 3 elements:
    1) <button id="ember359" type="button" data-easy-apply-next-button="" aria-label="Continue to next step" data-live-test-easy-apply-next-button="" class="artdeco-button artdeco-button--2 artdeco-button--primary ember-view">…</button> aka get_by_role("button", name="Continue to next step")
    2) <button id="ember270" type="button" aria-label="View next page" class="artdeco-button artdeco-button--muted artdeco-button--icon-right artdeco-button--1 artdeco-button--tertiary ember-view jobs-search-pagination__button jobs-search-pagination__button--next">…</button> aka get_by_label("View next page")
    3) <button id="ember313" type="button" aria-label="Company photos Next" data-control-name="COMPANY_LIFE_COMPANY_PHOTOS_NEXT" class="artdeco-button artdeco-button--circle artdeco-button--muted artdeco-button--icon-right artdeco-button--1 artdeco-button--tertiary ember-view artdeco-pagination__button artdeco-pagination__button--next">…</button> aka get_by_label("Company photos Next")
Call log:
  - waiting for locator("button:has-text('Next')")

'@
$ws.Range("F2").Value = $f2
$ws.Range("H2").Value = "19:34"

# Row 3 - strict mode violation error (2 elements), new time
$f3 = @'
Locator.click: Error: strict mode violation: locator("button:has-text('Next')") resolved to 2 elements:
    1) <button id="ember375" type="button" data-easy-apply-next-button="" aria-label="Continue to next step" data-live-test-easy-apply-next-button="" class="artdeco-button artdeco-button--2 artdeco-button--primary ember-view">…</button> aka get_by_role("button", name="Continue to next step")
    2) <button id="ember270" type="button" aria-label="View next page" class="artdeco-button artdeco-button--muted artdeco-button--icon-right artdeco-button--1 artdeco-button--tertiary ember-view jobs-search-pagination__button jobs-search-pagination__button--next">…</button> aka get_by_label("View next page")
Call log:
  - waiting for locator("button:has-text('Next')")

'@
$ws.Range("F3").Value = $f3
$ws.Range("H3").Value = "19:35"

# Row 4 - now successfully applied
$ws.Range("E4").Value = "Applied"
$ws.Range("F4").Value = "Success"
$ws.Range("H4").Value = "19:43"

# Row 5 - now successfully applied
$ws.Range("E5").Value = "Applied"
$ws.Range("F5").Value = "Success"
$ws.Range("H5").Value = "19:44"
